# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages are ready / complete:
#   * "Status" changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears.
#   * The "Latest Target File" (F) and "Latest Handback File" (G) columns
#     are now populated (with hyperlinks) on the zh-cn and de-de sheets.
#   * The "Latest Handback DateTime" (H) column is stamped with the
#     real handback timestamp for each language.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update the "Status" text everywhere (Overview + both language sheets)
# ---------------------------------------------------------------------
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (F) and
#    "Latest Handback File" (G) with hyperlinks, and stamp the
#    "Latest Handback DateTime" (H) column.
# ---------------------------------------------------------------------
$mdName  = "6f3021c9-7a13-43dc-9ed8-93935cb93275.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/70f39277c8ac817f12396a03261627114df07fc0/e2e/6f3021c9-7a13-43dc-9ed8-93935cb93275.md"

$zhcnXlfName = "6f3021c9-7a13-43dc-9ed8-93935cb93275.ec89bef60032987ec7d9f825313fe9e7d75db401.zh-cn.xlf"
$zhcnXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/426df2638fa5e88ce27d5ce9973d8ec7d3c469e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6f3021c9-7a13-43dc-9ed8-93935cb93275.ec89bef60032987ec7d9f825313fe9e7d75db401.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $mdUrl, "", "", $mdName)
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnXlfUrl, "", "", $zhcnXlfName)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $mdUrl, "", "", $mdName)
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnXlfUrl, "", "", $zhcnXlfName)

$zhcn.Range("H2").Value = "2016-03-11 14:43:58"
$zhcn.Range("H3").Value = "2016-03-11 14:43:58"

# ---------------------------------------------------------------------
# 3. de-de sheet: fill in "Latest Target File" (F) and
#    "Latest Handback File" (G) with hyperlinks, and stamp the
#    "Latest Handback DateTime" (H) column.
# ---------------------------------------------------------------------
$dedeXlfName = "6f3021c9-7a13-43dc-9ed8-93935cb93275.ec89bef60032987ec7d9f825313fe9e7d75db401.de-de.xlf"
$dedeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e50fd3ca806dbc767a5ed641c155f44ddf2002ea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6f3021c9-7a13-43dc-9ed8-93935cb93275.ec89bef60032987ec7d9f825313fe9e7d75db401.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("F2"), $mdUrl, "", "", $mdName)
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeXlfUrl, "", "", $dedeXlfName)
$dede.Hyperlinks.Add($dede.Range("F3"), $mdUrl, "", "", $mdName)
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeXlfUrl, "", "", $dedeXlfName)

$dede.Range("H2").Value = "2016-03-11 14:44:05"
$dede.Range("H3").Value = "2016-03-11 14:44:05"

Write-Host "Handback report generated."
